# Fix of following issues: 1. Reset Cart doesn't reset the cart.
# 2. All Main menu buttons remain enabled after reset Cart.
#
# This re-applies the "Reset Cart" action to the invoice template:
#  - bump the invoice DATE and INVOICE # to the new values
#  - clear out the customer-specific fields that should have been
#    wiped by a cart reset (bill-to name/location/contact + the
#    salesperson-table customer name/contact columns)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DATE (F2) and INVOICE # (F3)
$ws.Range("F2").Value = "2021-06-21 "
# INVOICE # is stored as text ("1" -> "3"); prefix with an apostrophe so it
# is entered as text rather than being auto-converted to a number, matching
# the original cell's text storage.
$ws.Range("F3").Value = "'3"

# Reset Cart: clear the previously populated customer info
$ws.Range("A9").ClearContents()
$ws.Range("A11").ClearContents()
$ws.Range("A13").ClearContents()
$ws.Range("B16").ClearContents()
$ws.Range("D16").ClearContents()
